$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text columns A:L to be stored as text, matching the source data's
# inlineStr typing (otherwise numeric-looking / date-looking strings would
# be auto-coerced to numbers/dates like real Excel does).
$ws.Range("A64:L64").NumberFormat = "@"

$ws.Range("A64").Value = "6248"
$ws.Range("B64").Value = "6/27/2025"
$ws.Range("C64").Value = "AVELLANEDA AV. 2395"
$ws.Range("D64").Value = "7"
$ws.Range("E64").Value = "807817952"
$ws.Range("F64").Value = "AYKO"
$ws.Range("G64").Value = "Pendiente"
$ws.Range("H64").Value = "Poste podrido"
$ws.Range("I64").Value = "1"
$ws.Range("J64").Value = "Cambio"
$ws.Range("K64").Value = "Sin equipos"
$ws.Range("L64").Value = "Poste"

# Restore the default (unstyled) formatting so the new row matches the
# rest of the sheet's data rows, which carry no explicit cell style
# (values remain text once entered, this only clears the display format).
$ws.Range("A64:L64").ClearFormats()

# Coordinate columns stay numeric, like the rest of the sheet.
$ws.Range("M64").Value = -58.4643
$ws.Range("N64").Value = -34.623993
